$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new row at 605 (pushes old rows 605..736 down to 606..737) ---
$ws.Rows(605).Insert()

# Expand Table1 so it covers the newly added row at the bottom (A8:K736 -> A8:K737)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K737"))

# The freshly inserted row (605) picks up generic default formatting instead
# of the table's row style, so copy the formatting from the row below (606,
# which now holds the data that used to live in row 605) across to row 605.
$ws.Range("A606:K606").Copy() | Out-Null
$ws.Range("A605:K605").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the calculated-column formula in G605 (PasteSpecial(Formats) does
# not bring formulas along).
$ws.Range("G605").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# The old footer row (formerly 736) lost its table-relative formula syntax
# when it was pushed outside Table1's old boundary mid-edit; put it back now
# that the table has been resized to include it again as row 737.
$ws.Range("G737").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- Populate the new row 605: one day Sick Leave taken 10/19/2023 ---
$ws.Range("B605").Value = "SL(1-0-0)"
$ws.Range("H605").Value = 1

# K605 needs the same date-formatted style already used elsewhere in column K
# (e.g. K600); copy its format rather than assigning NumberFormat directly so
# the existing style index is reused instead of a new one being minted.
$ws.Range("K600").Copy() | Out-Null
$ws.Range("K605").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("K605").Value = 45218

# --- Row 606 (was row 605 before the insert): Special Privilege Leave entry ---
$ws.Range("B606").Value = "SP(2-0-0)"
$ws.Range("K606").Value = "12/20,26/2023"

# --- Two EARNED (column C) entries of 1.25 each that were previously blank ---
$ws.Range("C603").Value = 1.25
$ws.Range("C604").Value = 1.25

$wb.Application.CalculateFull()
